$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.003.64'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.823.76'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Formula = "'311.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").Formula = "'0.4692"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("D8").Formula = "'0.3662"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("D9").Formula = "'0.07357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").Formula = "'0.8745"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").Formula = "'20.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").Value = '1.845.60'
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("D13").Formula = "'0.07310"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.27%  '
$ws.Range("D14").Formula = "'5.432"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Formula = "'6.524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").Formula = "'91.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Formula = "'0.000008747"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Formula = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '27.018.65'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("D24").Value = '2.065.70'
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").Formula = "'151.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").Formula = "'18.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Formula = "'2.144"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").Formula = "'116.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").Formula = "'0.7555"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").Formula = "'4.511"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").Formula = "'2.928"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").Formula = "'1.097"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").Formula = "'0.05313"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").Formula = "'0.01950"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Formula = "'2.979"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").Formula = "'2.373"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("D43").Formula = "'0.5309"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").Formula = "'8.486"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Formula = "'0.4897"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("D47").Formula = "'10.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").Formula = "'103.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").Formula = "'0.06299"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.52%  '
